$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns being updated, in order: B, C, D, E, F, I, K, N
$colLetters = @("B", "C", "D", "E", "F", "I", "K", "N")

# Each entry: row number followed by the 8 new values for the columns above
$data = @(
    @(2, 11.63507054624237, 7.930538139316975, 5.896791424070731, 16.35060966286245, 31.76627013736663, 24.76664481412141, 12.11681565905014, 19.76764990754434),
    @(3, 11.32922955996964, 7.622261269460787, 5.908790502108642, 15.4310494358278, 31.5921783426326, 24.76072662657355, 11.89048108504712, 19.82133918496889),
    @(4, 11.14079757440249, 7.429163852767203, 5.916309098691178, 14.84289575365247, 31.49368683623414, 24.76249397089074, 11.75299025392141, 19.85620113457738),
    @(5, 11.0639782084984, 7.349654403406963, 5.919411473705845, 14.59757357288768, 31.45569032339747, 24.76456984480222, 11.69741583384373, 19.87088472241143),
    @(6, 11.05122439681542, 7.336406663744019, 5.919928959359617, 14.55650614715822, 31.4495109729811, 24.76499633333522, 11.68821758029492, 19.87335173953777),
    @(7, 11.13976151449114, 7.428094690559357, 5.916350781975273, 14.83960970666786, 31.49316570618995, 24.76251648180608, 11.75223881016106, 19.85639723081888),
    @(8, 11.52982676912375, 7.825114753154973, 5.900897637999736, 16.03858008405622, 31.704516777151, 24.76348193371857, 12.03851762694498, 19.78576844610441),
    @(9, 12.2842034851214, 8.568011937093614, 5.871770276236965, 18.23143819052419, 32.18429079094873, 24.80831289210489, 12.60806915742565, 19.66230567844642),
    @(10, 12.82523132443722, 9.085787506507932, 5.851054636582234, 19.85965522486206, 32.57457172557507, 24.86749567343281, 13.02680132400796, 19.58075486323613),
    @(11, 13.06720586320511, 9.314169961490361, 5.841771971579222, 20.55897514653831, 32.75983199007965, 24.90011683907998, 13.21640717235771, 19.5456417708721),
    @(12, 13.15814373377796, 9.399548520056458, 5.838276587964337, 20.81787895845543, 32.83104958979806, 24.91328740250841, 13.28800505306839, 19.53263075300237),
    @(13, 13.13859089982853, 9.381210961838553, 5.839028510723496, 20.76238184345391, 32.81566506756843, 24.91041455932957, 13.27259530404959, 19.53542020608932),
    @(14, 13.07470178484429, 9.321216649589473, 5.841484011196473, 20.58039363580254, 32.76567006392588, 24.90118401603769, 13.22230203576053, 19.54456562288526),
    @(15, 13.03547491474057, 9.284322430181147, 5.842990634199595, 20.4681515212462, 32.73518376861094, 24.89563645370028, 13.19146750542517, 19.55020464620204),
    @(16, 12.80932589292058, 9.070711453911429, 5.851664073733452, 19.81312363710416, 32.56261594584329, 24.86547833813981, 13.01438593577018, 19.58308953642506),
    @(17, 12.66945833421396, 8.937777874002743, 5.857020693160667, 19.40071560867555, 32.45869598927793, 24.84843583829751, 12.90547246716317, 19.60377179033018),
    @(18, 12.58862565573061, 8.860648260627354, 5.860114979893014, 19.15961396647398, 32.39965351607953, 24.83916998708557, 12.8427500235165, 19.61585448263216),
    @(19, 12.56119414316188, 8.834420952439938, 5.861164952218847, 19.07731098687165, 32.37978942634491, 24.83612492252385, 12.82150219608573, 19.61997755514871),
    @(20, 12.68438800407222, 8.951998777527038, 5.856449099054441, 19.44502005921335, 32.4696832715827, 24.8501945220003, 12.9170751543624, 19.60155079453316),
    @(21, 13.09348709746862, 9.338868974866127, 5.840762239004252, 20.63400819840498, 32.7803263254451, 24.90387307861338, 13.23708044656714, 19.54187164131739),
    @(22, 13.35678052385554, 9.585242468286841, 5.830624905838215, 21.37662842324275, 32.98952565968997, 24.94371948928115, 13.44501415903945, 19.50453234551667),
    @(23, 13.21665882075076, 9.454362652382173, 5.83602504504814, 20.98341901854121, 32.87732284845578, 24.92201758795346, 13.33417056200223, 19.52430867593997),
    @(24, 12.67763960402965, 8.945571695017831, 5.856707470905611, 19.42500248796413, 32.46471372947023, 24.84939776343452, 12.91182990582768, 19.60255430804283),
    @(25, 12.08198253152991, 8.371555263831089, 5.879527460027903, 17.63262602567066, 32.04770359443615, 24.79157776242118, 12.45360935518264, 19.6940967928718)
)

foreach ($entry in $data) {
    $rowNum = $entry[0]
    for ($i = 0; $i -lt $colLetters.Length; $i++) {
        $cellRef = "$($colLetters[$i])$rowNum"
        $ws.Range($cellRef).Value = $entry[$i + 1]
    }
}
